$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded, so insert a row for it right
# after the header/most-recent row (row 3), shifting the existing rows
# (old 3-13) down to (4-14). Excel's Insert copies the row-above formatting
# automatically (matches the date-style cell D2 -> D3).
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new observation's data. The
# non-numeric/date columns mirror the rest of the table for this market.
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 'Macroferia Regional de Talca'
$ws.Range("C3").Value = 'Maule'
$ws.Range("D3").Value = 44453
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 'Fruta'
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = 'Otros'
$ws.Range("I3").Value = 100107002
$ws.Range("J3").Value = 'Chirimoya'
$ws.Range("K3").Value = 'Cultivar IV Región'
$ws.Range("L3").Value = 'Especial'
$ws.Range("M3").Value = 135
$ws.Range("N3").Value = 30000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 30000
$ws.Range("Q3").Value = '$/bandeja 10 kilos'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 3000
$ws.Range("T3").Value = 10
